$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Georgia (row 5) is now completed: set patient count and mark it green (completed fill)
$ws.Range("B5").Value = 3100
$ws.Range("B5").Interior.Color = 5287936   # RGB(0,176,80) -> matches the "completed" green fill

# Ukraine (row 14) moves from "not started" (plain) to "in progress" (orange fill)
$ws.Range("B14").Interior.Color = 49407    # RGB(255,192,0) -> matches the "in progress" orange fill

# Update the last selected cell to reflect where the user finished editing
$ws.Range("D14").Select()
